$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "28.213.15"
Set-TextValue "E2" "  -1.38%  "
Set-TextValue "D3" "1.841.10"
Set-TextValue "E3" "  -0.46%  "
Set-TextValue "D4" "1.006"
Set-TextValue "E4" "  +0.35%  "
Set-TextValue "D5" "326.27"
Set-TextValue "E5" "  -2.81%  "
Set-TextValue "D6" "1.006"
Set-TextValue "E6" "  +0.45%  "
Set-TextValue "D7" "0.4649"
Set-TextValue "E7" "  -0.08%  "
Set-TextValue "E8" "  -0.80%  "
Set-TextValue "D9" "0.07863"
Set-TextValue "E9" "  -0.49%  "
Set-TextValue "D10" "0.9646"
Set-TextValue "E10" "  -1.37%  "
Set-TextValue "D11" "22.14"
Set-TextValue "E11" "  -0.69%  "
Set-TextValue "B12" "Polkadot"
Set-TextValue "C12" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D12" "5.704"
Set-TextValue "E12" "  -1.84%  "
Set-TextValue "B13" "Chainlink"
Set-TextValue "C13" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D13" "6.878"
Set-TextValue "E13" "  -1.22%  "
Set-TextValue "B14" "WrappedEther"
Set-TextValue "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.749.17"
Set-TextValue "E14" "  -7.74%  "
Set-TextValue "D15" "0.06912"
Set-TextValue "E15" "  +0.15%  "
Set-TextValue "D16" "88.66"
Set-TextValue "E16" "  +1.05%  "
Set-TextValue "D17" "1.006"
Set-TextValue "E17" "  +0.40%  "
Set-TextValue "D18" "0.000009971"
Set-TextValue "D19" "16.76"
Set-TextValue "E19" "  -1.67%  "
Set-TextValue "D20" "1.005"
Set-TextValue "E20" "  +0.28%  "
Set-TextValue "D21" "28.215.26"
Set-TextValue "E21" "  -1.44%  "
Set-TextValue "D22" "5.303"
Set-TextValue "E22" "  -1.47%  "
Set-TextValue "D23" "11.05"
Set-TextValue "E23" "  -1.50%  "
Set-TextValue "D24" "2.106"
Set-TextValue "E24" "  -2.15%  "
Set-TextValue "D25" "2.077.78"
Set-TextValue "E25" "  -1.62%  "
Set-TextValue "D26" "154.54"
Set-TextValue "E26" "  +0.94%  "
Set-TextValue "D27" "19.18"
Set-TextValue "E27" "  -0.93%  "
Set-TextValue "D28" "5.755"
Set-TextValue "E28" "  -4.92%  "
Set-TextValue "D29" "1.968"
Set-TextValue "E29" "  -2.11%  "
Set-TextValue "D30" "118.95"
Set-TextValue "E30" "  +1.32%  "
Set-TextValue "D31" "0.09264"
Set-TextValue "E31" "  -0.98%  "
Set-TextValue "D32" "0.9328"
Set-TextValue "E32" "  -3.58%  "
Set-TextValue "D33" "5.287"
Set-TextValue "E33" "  -1.41%  "
Set-TextValue "D34" "1.328"
Set-TextValue "E34" "  -1.42%  "
Set-TextValue "D35" "3.335"
Set-TextValue "E35" "  -3.74%  "
Set-TextValue "D36" "0.05822"
Set-TextValue "E36" "  -4.63%  "
Set-TextValue "D37" "0.02121"
Set-TextValue "E37" "  -3.61%  "
Set-TextValue "D38" "1.139"
Set-TextValue "E38" "  -1.94%  "
Set-TextValue "D39" "7.787"
Set-TextValue "E39" "  +1.47%  "
Set-TextValue "D40" "0.5593"
Set-TextValue "E40" "  -1.82%  "
Set-TextValue "D41" "9.927"
Set-TextValue "E41" "  -1.65%  "
Set-TextValue "D42" "0.1765"
Set-TextValue "E42" "  -1.52%  "
Set-TextValue "D43" "0.07272"
Set-TextValue "E43" "  +2.56%  "
Set-TextValue "D44" "11.61"
Set-TextValue "E44" "  -1.08%  "
Set-TextValue "D45" "0.5280"
Set-TextValue "E45" "  -1.87%  "
Set-TextValue "E46" "  -8.74%  "
Set-TextValue "E47" "  -11.88%  "
Set-TextValue "D48" "1.838"
Set-TextValue "E48" "  -3.35%  "
Set-TextValue "D49" "114.09"
Set-TextValue "E49" "  +0.87%  "
Set-TextValue "D50" "1.006"
Set-TextValue "E50" "  +0.60%  "
Set-TextValue "D51" "2.327"
Set-TextValue "E51" "  -0.71%  "
